$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("E2").Value = 2
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 1.902492
$ws.Range("H2").Value = 3.804984
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.5
$ws.Range("M2").Value = 0.09207700000000001
$ws.Range("N2").Value = 0.184154
$ws.Range("O2").Value = 0.0789959771480734
$ws.Range("P2").Value = 0.05545240531440215
$ws.Range("Q2").Value = 0.175175755884
$ws.Range("R2").Value = 0.7007030235360001
$ws.Range("S2").Value = 0.0789959771480734
$ws.Range("T2").Value = 0.05545240531440215

# Row 3 updates
$ws.Range("E3").Value = 2
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 1.902492
$ws.Range("H3").Value = 3.804984
$ws.Range("O3").Value = 0.8491451975864605
$ws.Range("P3").Value = 0.8941052196698643
$ws.Range("Q3").Value = 1.883002871952
$ws.Range("R3").Value = 11.298017231712
$ws.Range("S3").Value = 0.8491451975864605
$ws.Range("T3").Value = 0.8941052196698643

# Row 4 updates
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 1.902492
$ws.Range("H4").Value = 3.804984
$ws.Range("M4").Value = 0.083758
$ws.Range("N4").Value = 0.167516
$ws.Range("O4").Value = 0.07185882526546619
$ws.Range("P4").Value = 0.05044237501573352
$ws.Range("Q4").Value = 0.159348924936
$ws.Range("R4").Value = 0.637395699744
$ws.Range("S4").Value = 0.07185882526546619
$ws.Range("T4").Value = 0.05044237501573352

# Remove row 5 entirely (MuSCs/Fgf5/Fgfr2/Neutrophils pair no longer present)
$ws.Rows("5:5").Delete()
